$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = [double]"0.0001000210781354453"
$ws.Range("E3").Value = [double]"0.0001000210781354453"

$ws.Range("D4").Value = [double]"0.0004044725105451163"
$ws.Range("E4").Value = [double]"0.0004044725105451163"

$ws.Range("D5").Value = [double]"1.175778944526739E-27"
$ws.Range("E5").Value = [double]"1.175778944526739E-27"

$ws.Range("D6").Value = [double]"4.368558616677276E-26"
$ws.Range("E6").Value = [double]"4.368558616677276E-26"

$ws.Range("D7").Value = [double]"0.9999999999955131"
$ws.Range("E7").Value = [double]"4.486855331720108E-12"

$ws.Range("D8").Value = [double]"0.995927927319859"
$ws.Range("E8").Value = [double]"0.004072072680140981"

$ws.Range("D9").Value = [double]"0.9999999166203365"
$ws.Range("E9").Value = [double]"8.337966350691062E-08"

$ws.Range("D10").Value = [double]"7.183843074994926E-17"
$ws.Range("E10").Value = [double]"0.9999999999999999"

$ws.Range("D11").Value = [double]"1.104957013213112E-60"
$ws.Range("F11").Value = [double]"24.39967727661133"
